$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete row 26 (RM 232) entirely - rows below shift up
$ws.Rows("26").Delete()

# 2) Delete what is now row 27 (originally SC 92) entirely - rows below shift up
$ws.Rows("27").Delete()

# 3) Apply individual cell edits for column D (and some column C) per the diff

# Row 6: D6 becomes numeric -14.2
$ws.Range("D6").Value = -14.2

# Row 8: D8 becomes blank (was numeric -13.9)
$ws.Range("D8").ClearContents()

# Row 12: D12 becomes numeric -14.1
$ws.Range("D12").Value = -14.1

# Row 14: D14 becomes blank (was numeric -13.1)
$ws.Range("D14").ClearContents()

# Row 17: D17 becomes numeric -14.7
$ws.Range("D17").Value = -14.7

# Row 18: D18 becomes numeric -15.2
$ws.Range("D18").Value = -15.2

# Row 19: D19 becomes blank (was numeric -15.5)
$ws.Range("D19").ClearContents()

# Row 20: D20 becomes blank (was numeric -14)
$ws.Range("D20").ClearContents()

# Row 23: D23 becomes numeric -13.9
$ws.Range("D23").Value = -13.9

# Rows 26-33 (post row-deletion, shifted up) additional edits:
# Row 27 (SC 101): C27 becomes numeric 10; D27 becomes blank (was -14.6)
$ws.Range("C27").Value = 10
$ws.Range("D27").ClearContents()

# Row 28 (SC 105): C28 becomes blank (was 11.1)
$ws.Range("C28").ClearContents()

# Row 29 (SC 119): C29 becomes blank (was 11.2)
$ws.Range("C29").ClearContents()

# Row 30 (SC 120): C30 becomes numeric 11.4
$ws.Range("C30").Value = 11.4

# Row 32 (SC 193): C32 becomes blank (was 10.5)
$ws.Range("C32").ClearContents()
